# Add data for 2021-12-14: the running "through December 05" totals are
# updated to "through December 06", picking up one additional day's worth
# of carjacking incidents across several neighborhoods/months.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date.
$ws.Name = "Through 2021-12-06"

# Update the header label (column B, row 1) to match.
$ws.Range("B1").Value = "December 2021 (through December 06)"

# West Town: new incident in December 2017 (column AX)
$ws.Range("AX2").Value = 1

# Englewood: new incident in current month, December 2021 (column B)
$ws.Range("B3").Value = 1

# North Lawndale: December 2020 (column N) count revised upward
$ws.Range("N4").Value = 4

# Grand Crossing: new incident in December 2020 (column N)
$ws.Range("N8").Value = 1

# Humboldt Park: December 2017 (column AX) count revised upward
$ws.Range("AX9").Value = 2

# Chatham: new incident in December 2019 (column Z)
$ws.Range("Z13").Value = 1

# Chicago Lawn: new incident in current month, December 2021 (column B)
$ws.Range("B14").Value = 1

# United Center: new incident in December 2017 (column AX)
$ws.Range("AX16").Value = 1

# Bridgeport: new incident in current month, December 2021 (column B)
$ws.Range("B17").Value = 1

# Grand Boulevard: new incidents in current month, December 2021 (column B)
$ws.Range("B18").Value = 3

# Little Village: new incident in December 2017 (column AX)
$ws.Range("AX20").Value = 1

# South Shore: current month (B) revised upward, plus new incidents in
# December 2020 (N) and December 2016 (BJ)
$ws.Range("B21").Value = 3
$ws.Range("N21").Value = 1
$ws.Range("BJ21").Value = 2

# South Chicago: new incident in December 2020 (column N)
$ws.Range("N32").Value = 1

# East Village: new incident in current month, December 2021 (column B)
$ws.Range("B41").Value = 1

# Riverdale: current month (column B) revised upward
$ws.Range("B48").Value = 2

# Avondale: current month (column B) revised upward
$ws.Range("B67").Value = 4

# Belmont Cragin: new incident in December 2016 (column BJ)
$ws.Range("BJ68").Value = 1

# Printers Row: new incident in current month, December 2021 (column B)
$ws.Range("B92").Value = 1
